$wb = $excel.ActiveWorkbook

# The single "Sheet1" becomes "Examinees"; the title row is dropped and the
# header row (former row 2) moves up to row 1, gaining two new columns.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Examinees"

# Drop the "1st Dan Candidates" title row so the header row becomes row 1.
$ws1.Rows.Item(1).Delete() | Out-Null

# Append the two new header columns.
$ws1.Range("H1").Value = "Group"
$ws1.Range("I1").Value = "New Rank"

# Match the recorded selection on the Examinees sheet.
$ws1.Range("I2").Select() | Out-Null

# Add the new, still-empty "Examiners" sheet right after "Examinees".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Examiners"

# Leave focus back on the Examinees sheet/tab.
$ws1.Activate() | Out-Null
